$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Isosteric Heat Analysis")

$ws.Cells.Item(19, 2).Value = -18378.575811054645
$ws.Cells.Item(19, 3).Value = 579.5134118970051
$ws.Cells.Item(20, 2).Value = -18514.235348989187
$ws.Cells.Item(20, 3).Value = 496.3133064403586
$ws.Cells.Item(21, 2).Value = -18657.396694491377
$ws.Cells.Item(21, 3).Value = 408.4859772109164
$ws.Cells.Item(22, 2).Value = -18808.52702744859
$ws.Cells.Item(22, 3).Value = 316.3657643973871
$ws.Cells.Item(23, 2).Value = -18968.09808708953
$ws.Cells.Item(23, 3).Value = 221.76317533168935
$ws.Cells.Item(24, 2).Value = -19136.574085987635
$ws.Cells.Item(24, 3).Value = 133.8686923887531
$ws.Cells.Item(25, 2).Value = -19314.394300115335
$ws.Cells.Item(25, 3).Value = 103.9677643832511
$ws.Cells.Item(26, 2).Value = -19501.94852921517
$ws.Cells.Item(26, 3).Value = 182.0824505092694
$ws.Cells.Item(27, 2).Value = -19699.543131485745
$ws.Cells.Item(27, 3).Value = 301.9546068240059
$ws.Cells.Item(28, 2).Value = -19907.35478495998
$ws.Cells.Item(28, 3).Value = 438.28138925260197
$ws.Cells.Item(29, 2).Value = -20125.368576439698
$ws.Cells.Item(29, 3).Value = 586.6344217497052
$ws.Cells.Item(30, 2).Value = -20353.296595534208
$ws.Cells.Item(30, 3).Value = 746.0290239942018
$ws.Cells.Item(31, 2).Value = -20590.473158177017
$ws.Cells.Item(31, 3).Value = 916.1892658911182
$ws.Cells.Item(32, 2).Value = -20835.72351462588
$ws.Cells.Item(32, 3).Value = 1096.8900903129013
$ws.Cells.Item(33, 2).Value = -21087.20505578697
$ws.Cells.Item(33, 3).Value = 1287.6808536999324
$ws.Cells.Item(34, 2).Value = -21342.224515142654
$ws.Cells.Item(34, 3).Value = 1487.6929961186208
$ws.Cells.Item(35, 2).Value = -21597.042519167004
$ws.Cells.Item(35, 3).Value = 1695.466412031768
$ws.Cells.Item(36, 2).Value = -21846.688900679317
$ws.Cells.Item(36, 3).Value = 1908.790186204067
$ws.Cells.Item(37, 2).Value = -22084.828273975676
$ws.Cells.Item(37, 3).Value = 2124.583146979235
$ws.Cells.Item(38, 2).Value = -22303.732977775417
$ws.Cells.Item(38, 3).Value = 2338.8626363762783
$ws.Cells.Item(39, 2).Value = -22494.433376419212
$ws.Cells.Item(39, 3).Value = 2546.865510271392
$ws.Cells.Item(40, 2).Value = -22647.113308454293
$ws.Cells.Item(40, 3).Value = 2743.382954097288
$ws.Cells.Item(41, 2).Value = -22751.789332519074
$ws.Cells.Item(41, 3).Value = 2923.3360620471194
$ws.Cells.Item(42, 2).Value = -22799.249974716393
$ws.Cells.Item(42, 3).Value = 3082.5454355538914
$ws.Cells.Item(43, 2).Value = -22782.143891614673
$ws.Cells.Item(43, 3).Value = 3218.549326920287
